$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 2138.889  # H17
$ws.Cells.Item(17, 10).Value = 2121.4285  # J17
$ws.Cells.Item(17, 12).Value = 6364.2855  # L17
$ws.Cells.Item(17, 14).Value = -6700.2855  # N17

$ws.Cells.Item(18, 8).Value = 3599.8  # H18
$ws.Cells.Item(18, 9).Value = 2749.75  # I18
$ws.Cells.Item(18, 10).Value = 7000  # J18
$ws.Cells.Item(18, 11).Value = 2749.75  # K18
$ws.Cells.Item(18, 12).Value = 7000  # L18
$ws.Cells.Item(18, 13).Value = -2465.75  # M18
$ws.Cells.Item(18, 14).Value = -7568  # N18

$ws.Cells.Item(53, 8).Value = 156.125  # H53
$ws.Cells.Item(53, 10).Value = 237.5  # J53
$ws.Cells.Item(53, 12).Value = 237.5  # L53
$ws.Cells.Item(53, 14).Value = -1511.5  # N53

$ws.Cells.Item(92, 8).Value = 1401.7646  # H92
$ws.Cells.Item(92, 9).Value = 1411.9  # I92
$ws.Cells.Item(92, 11).Value = 1411.9  # K92
$ws.Cells.Item(92, 13).Value = -163.9000000000001  # M92

$ws.Cells.Item(96, 8).Value = 2279.9  # H96
$ws.Cells.Item(96, 9).Value = 2518.1667  # I96
$ws.Cells.Item(96, 10).Value = 1922.5  # J96
$ws.Cells.Item(96, 11).Value = 7554.500100000001  # K96
$ws.Cells.Item(96, 12).Value = 5767.5  # L96
$ws.Cells.Item(96, 13).Value = -6181.500100000001  # M96
$ws.Cells.Item(96, 14).Value = -8513.5  # N96

$ws.Cells.Item(103, 8).Value = 1300  # H103
$ws.Cells.Item(103, 9).Value = 1000  # I103
$ws.Cells.Item(103, 10).Value = 2500  # J103
$ws.Cells.Item(103, 11).Value = 3000  # K103
$ws.Cells.Item(103, 12).Value = 7500  # L103
$ws.Cells.Item(103, 13).Value = -2414  # M103
$ws.Cells.Item(103, 14).Value = -8672  # N103

$ws.Cells.Item(137, 8).Value = 1698  # H137
$ws.Cells.Item(137, 9).Value = 1553.625  # I137
$ws.Cells.Item(137, 10).Value = 2083  # J137
$ws.Cells.Item(137, 11).Value = 4660.875  # K137
$ws.Cells.Item(137, 12).Value = 6249  # L137
$ws.Cells.Item(137, 13).Value = -2110.875  # M137
$ws.Cells.Item(137, 14).Value = -11349  # N137

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5634.778  # H32
$ws.Cells.Item(32, 9).Value = 5634.778  # I32
$ws.Cells.Item(32, 11).Value = 5634.778  # K32
$ws.Cells.Item(32, 13).Value = -5347.778  # M32

$ws.Cells.Item(97, 8).Value = 2031.3334  # H97
$ws.Cells.Item(97, 9).Value = 872  # I97
$ws.Cells.Item(97, 11).Value = 872  # K97
$ws.Cells.Item(97, 13).Value = -376  # M97

$ws.Cells.Item(106, 8).Value = 0  # H106
$ws.Cells.Item(106, 10).Value = 0  # J106
$ws.Cells.Item(106, 12).Value = 0  # L106
$ws.Cells.Item(106, 14).ClearContents()  # N106

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(35, 8).Value = 15000  # H35
$ws.Cells.Item(35, 9).Value = 0  # I35
$ws.Cells.Item(35, 10).Value = 15000  # J35
$ws.Cells.Item(35, 11).Value = 0  # K35
$ws.Cells.Item(35, 12).Value = 15000  # L35
$ws.Cells.Item(35, 13).ClearContents()  # M35
$ws.Cells.Item(35, 14).Value = -15620  # N35

$ws.Cells.Item(36, 8).Value = 0  # H36
$ws.Cells.Item(36, 9).Value = 0  # I36
$ws.Cells.Item(36, 11).Value = 0  # K36
$ws.Cells.Item(36, 13).ClearContents()  # M36

$ws.Cells.Item(86, 8).Value = 1394  # H86
$ws.Cells.Item(86, 9).Value = 1394  # I86
$ws.Cells.Item(86, 11).Value = 1394  # K86
$ws.Cells.Item(86, 13).Value = -271  # M86

$ws.Cells.Item(89, 8).Value = 1394  # H89
$ws.Cells.Item(89, 9).Value = 1394  # I89
$ws.Cells.Item(89, 11).Value = 6970  # K89
$ws.Cells.Item(89, 13).Value = -1354  # M89

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(12, 8).Value = 930  # H12
$ws.Cells.Item(12, 9).Value = 645  # I12
$ws.Cells.Item(12, 10).Value = 1500  # J12
$ws.Cells.Item(12, 11).Value = 645  # K12
$ws.Cells.Item(12, 12).Value = 1500  # L12
$ws.Cells.Item(12, 13).Value = -475  # M12
$ws.Cells.Item(12, 14).Value = -1840  # N12

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 633163.5  # H4
$ws.Cells.Item(4, 9).Value = 703507.25  # I4
$ws.Cells.Item(4, 11).Value = 2110521.75  # K4
$ws.Cells.Item(4, 13).Value = -2110409.75  # M4

$ws.Cells.Item(22, 8).Value = 1000  # H22
$ws.Cells.Item(22, 9).Value = 1000  # I22
$ws.Cells.Item(22, 11).Value = 3000  # K22
$ws.Cells.Item(22, 13).Value = -2831  # M22

$ws.Cells.Item(27, 8).Value = 1000  # H27
$ws.Cells.Item(27, 9).Value = 1000  # I27
$ws.Cells.Item(27, 11).Value = 3000  # K27
$ws.Cells.Item(27, 13).Value = -2898  # M27

$ws.Cells.Item(32, 8).Value = 1000  # H32
$ws.Cells.Item(32, 10).Value = 1000  # J32
$ws.Cells.Item(32, 12).Value = 3000  # L32
$ws.Cells.Item(32, 14).Value = -3566  # N32

$ws.Cells.Item(39, 8).Value = 5000  # H39
$ws.Cells.Item(39, 10).Value = 5000  # J39
$ws.Cells.Item(39, 12).Value = 15000  # L39
$ws.Cells.Item(39, 14).Value = -15588  # N39

$ws.Cells.Item(40, 8).Value = 203.77777  # H40
$ws.Cells.Item(40, 9).Value = 80  # I40
$ws.Cells.Item(40, 11).Value = 320  # K40
$ws.Cells.Item(40, 13).Value = -251  # M40

$ws.Cells.Item(55, 8).Value = 2039.45  # H55
$ws.Cells.Item(55, 10).Value = 2700  # J55
$ws.Cells.Item(55, 12).Value = 8100  # L55
$ws.Cells.Item(55, 14).Value = -8454  # N55

$ws.Cells.Item(95, 8).Value = 1000  # H95
$ws.Cells.Item(95, 9).Value = 1000  # I95
$ws.Cells.Item(95, 11).Value = 3000  # K95
$ws.Cells.Item(95, 13).Value = -941  # M95

$ws.Cells.Item(139, 8).Value = 3375  # H139
$ws.Cells.Item(139, 9).Value = 2562.5  # I139
$ws.Cells.Item(139, 10).Value = 5000  # J139
$ws.Cells.Item(139, 11).Value = 7687.5  # K139
$ws.Cells.Item(139, 12).Value = 15000  # L139
$ws.Cells.Item(139, 13).Value = -2547.5  # M139
$ws.Cells.Item(139, 14).Value = -25280  # N139

$ws.Cells.Item(141, 8).Value = 2606.8  # H141
$ws.Cells.Item(141, 9).Value = 2606.8  # I141
$ws.Cells.Item(141, 11).Value = 7820.400000000001  # K141
$ws.Cells.Item(141, 13).Value = -2640.400000000001  # M141

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 466.66666  # H113
$ws.Cells.Item(113, 9).Value = 450  # I113
$ws.Cells.Item(113, 11).Value = 450  # K113
$ws.Cells.Item(113, 13).Value = 1720  # M113

$ws.Cells.Item(122, 8).Value = 4199.7144  # H122
$ws.Cells.Item(122, 9).Value = 4566.3335  # I122
$ws.Cells.Item(122, 11).Value = 13699.0005  # K122
$ws.Cells.Item(122, 13).Value = -11249.0005  # M122

$ws.Cells.Item(132, 8).Value = 4254.222  # H132
$ws.Cells.Item(132, 10).Value = 4381.3335  # J132
$ws.Cells.Item(132, 12).Value = 13144.0005  # L132
$ws.Cells.Item(132, 14).Value = -18204.0005  # N132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 999  # H7
$ws.Cells.Item(7, 9).Value = 999  # I7
$ws.Cells.Item(7, 11).Value = 999  # K7
$ws.Cells.Item(7, 13).Value = -887  # M7

$ws.Cells.Item(22, 8).Value = 1004.4286  # H22
$ws.Cells.Item(22, 9).Value = 963.5  # I22
$ws.Cells.Item(22, 10).Value = 1250  # J22
$ws.Cells.Item(22, 11).Value = 963.5  # K22
$ws.Cells.Item(22, 12).Value = 1250  # L22
$ws.Cells.Item(22, 13).Value = -668.5  # M22
$ws.Cells.Item(22, 14).Value = -1840  # N22

$ws.Cells.Item(27, 8).Value = 1004.4286  # H27
$ws.Cells.Item(27, 9).Value = 963.5  # I27
$ws.Cells.Item(27, 10).Value = 1250  # J27
$ws.Cells.Item(27, 11).Value = 963.5  # K27
$ws.Cells.Item(27, 12).Value = 1250  # L27
$ws.Cells.Item(27, 13).Value = -856.5  # M27
$ws.Cells.Item(27, 14).Value = -1464  # N27

$ws.Cells.Item(61, 8).Value = 3669.8  # H61
$ws.Cells.Item(61, 9).Value = 3987.25  # I61
$ws.Cells.Item(61, 10).Value = 2400  # J61
$ws.Cells.Item(61, 11).Value = 3987.25  # K61
$ws.Cells.Item(61, 12).Value = 2400  # L61
$ws.Cells.Item(61, 13).Value = -3785.25  # M61
$ws.Cells.Item(61, 14).Value = -2804  # N61

$ws.Cells.Item(95, 8).Value = 20000  # H95
$ws.Cells.Item(95, 10).Value = 20000  # J95
$ws.Cells.Item(95, 12).Value = 20000  # L95
$ws.Cells.Item(95, 14).Value = -25492  # N95

$ws.Cells.Item(113, 8).Value = 3669.8  # H113
$ws.Cells.Item(113, 9).Value = 3987.25  # I113
$ws.Cells.Item(113, 10).Value = 2400  # J113
$ws.Cells.Item(113, 11).Value = 3987.25  # K113
$ws.Cells.Item(113, 12).Value = 2400  # L113
$ws.Cells.Item(113, 13).Value = -1817.25  # M113
$ws.Cells.Item(113, 14).Value = -6740  # N113

$ws.Cells.Item(126, 8).Value = 999  # H126
$ws.Cells.Item(126, 9).Value = 999  # I126
$ws.Cells.Item(126, 11).Value = 2997  # K126
$ws.Cells.Item(126, 13).Value = -527  # M126

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(11, 8).Value = 3768  # H11
$ws.Cells.Item(11, 9).Value = 2752  # I11
$ws.Cells.Item(11, 10).Value = 5800  # J11
$ws.Cells.Item(11, 11).Value = 2752  # K11
$ws.Cells.Item(11, 12).Value = 5800  # L11
$ws.Cells.Item(11, 13).Value = -2610  # M11
$ws.Cells.Item(11, 14).Value = -6084  # N11

$ws.Cells.Item(30, 8).Value = 19504.5  # H30
$ws.Cells.Item(30, 9).Value = 16006  # I30
$ws.Cells.Item(30, 11).Value = 16006  # K30
$ws.Cells.Item(30, 13).Value = -15899  # M30

$ws.Cells.Item(31, 8).Value = 26339  # H31
$ws.Cells.Item(31, 9).Value = 24508.5  # I31
$ws.Cells.Item(31, 11).Value = 24508.5  # K31
$ws.Cells.Item(31, 13).Value = -24160.5  # M31

$ws.Cells.Item(55, 8).Value = 48  # H55
$ws.Cells.Item(55, 9).Value = 48  # I55
$ws.Cells.Item(55, 11).Value = 48  # K55
$ws.Cells.Item(55, 13).Value = 229  # M55

$ws.Cells.Item(99, 8).Value = 50000  # H99
$ws.Cells.Item(99, 9).Value = 50000  # I99
$ws.Cells.Item(99, 11).Value = 50000  # K99
$ws.Cells.Item(99, 13).Value = -47005  # M99

$ws.Cells.Item(100, 8).Value = 8162.4287  # H100
$ws.Cells.Item(100, 9).Value = 11297  # I100
$ws.Cells.Item(100, 10).Value = 2520.2  # J100
$ws.Cells.Item(100, 11).Value = 22594  # K100
$ws.Cells.Item(100, 12).Value = 5040.4  # L100
$ws.Cells.Item(100, 13).Value = -22053  # M100
$ws.Cells.Item(100, 14).Value = -6122.4  # N100

$ws.Cells.Item(113, 8).Value = 999.6667  # H113
$ws.Cells.Item(113, 9).Value = 1200  # I113
$ws.Cells.Item(113, 11).Value = 3600  # K113
$ws.Cells.Item(113, 13).Value = -1430  # M113

$ws.Cells.Item(126, 8).Value = 3174.75  # H126
$ws.Cells.Item(126, 9).Value = 2233  # I126
$ws.Cells.Item(126, 11).Value = 6699  # K126
$ws.Cells.Item(126, 13).Value = -4229  # M126
